$d = $word.ActiveDocument

# 1. Update the date text from 2024-04-08 to 2024-04-09
$d.Content.Find.Execute("2024-04-08", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-04-09", 2)

# 2. Remove the trailing image/source paragraphs that followed the
#    "...discussed." paragraph (pie_chart image, its source link,
#    the two bar_plot images, and their shared source link).
$start = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7).EndsWith("discussed.")) {
        $start = $i + 1
        break
    }
}

if ($start -ne $null -and $start -le $d.Paragraphs.Count) {
    $pStart = $d.Paragraphs.Item($start)
    $pEnd = $d.Paragraphs.Item($d.Paragraphs.Count)
    $r = $d.Range($pStart.Range.Start, $pEnd.Range.End)
    $r.Delete()
}
